$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 91/92 (Cartagines-Saprissa <-> Sporting San Jose-Guanacasteca) were
#    reordered upstream: swap everything except the shared index/date columns.
# ---------------------------------------------------------------------------
$row91 = $ws.Range("F91:V91").Value()
$row92 = $ws.Range("F92:V92").Value()
$ws.Range("F91:V91").Value = $row92
$ws.Range("F92:V92").Value = $row91

# ---------------------------------------------------------------------------
# 2) Rows 110/111 (AD Santos-San Carlos <-> Alajuelense-Grecia) likewise swap.
# ---------------------------------------------------------------------------
$row110 = $ws.Range("F110:V110").Value()
$row111 = $ws.Range("F111:V111").Value()
$ws.Range("F110:V110").Value = $row111
$ws.Range("F111:V111").Value = $row110

# ---------------------------------------------------------------------------
# 3) Two new matches were appended at the bottom of the sheet: rows 126/127.
#    Clone the formatting of the last existing data row (125) first so the
#    index column keeps its bold/border style and the date column keeps its
#    custom number format, then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A125:V125").Copy()
$ws.Range("A126:V126").PasteSpecial(-4122)
$ws.Range("A125:V125").Copy()
$ws.Range("A127:V127").PasteSpecial(-4122)

$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 2).Value = "costa-rica"
$ws.Cells.Item(126, 3).Value = "primera-division"
$ws.Cells.Item(126, 4).Value = "2023-2024"
$ws.Cells.Item(126, 5).Value = 45253.91666666666
$ws.Cells.Item(126, 6).Value = "Grecia"
$ws.Cells.Item(126, 7).Value = 1
$ws.Cells.Item(126, 8).Value = "Saprissa"
$ws.Cells.Item(126, 9).Value = 2
$ws.Cells.Item(126, 10).Value = 6.47
$ws.Cells.Item(126, 11).Value = "16/11/2023 22:12"
$ws.Cells.Item(126, 12).Value = 7.81
$ws.Cells.Item(126, 13).Value = "23/11/2023 21:58"
$ws.Cells.Item(126, 14).Value = 4.27
$ws.Cells.Item(126, 15).Value = "16/11/2023 22:12"
$ws.Cells.Item(126, 16).Value = 4.44
$ws.Cells.Item(126, 17).Value = "23/11/2023 21:58"
$ws.Cells.Item(126, 18).Value = 1.44
$ws.Cells.Item(126, 19).Value = "16/11/2023 22:12"
$ws.Cells.Item(126, 20).Value = 1.37
$ws.Cells.Item(126, 21).Value = "23/11/2023 21:29"
$ws.Cells.Item(126, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/grecia-saprissa/f7eU5K9q/"

$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = "costa-rica"
$ws.Cells.Item(127, 3).Value = "primera-division"
$ws.Cells.Item(127, 4).Value = "2023-2024"
$ws.Cells.Item(127, 5).Value = 45254.08333333334
$ws.Cells.Item(127, 6).Value = "San Carlos"
$ws.Cells.Item(127, 7).Value = 1
$ws.Cells.Item(127, 8).Value = "Zeledon"
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 10).Value = 1.55
$ws.Cells.Item(127, 11).Value = "17/11/2023 03:12"
$ws.Cells.Item(127, 12).Value = 1.3
$ws.Cells.Item(127, 13).Value = "24/11/2023 01:55"
$ws.Cells.Item(127, 14).Value = 4.21
$ws.Cells.Item(127, 15).Value = "17/11/2023 03:12"
$ws.Cells.Item(127, 16).Value = 5.75
$ws.Cells.Item(127, 17).Value = "24/11/2023 01:57"
$ws.Cells.Item(127, 18).Value = 5.09
$ws.Cells.Item(127, 19).Value = "17/11/2023 03:12"
$ws.Cells.Item(127, 20).Value = 9.16
$ws.Cells.Item(127, 21).Value = "24/11/2023 01:57"
$ws.Cells.Item(127, 22).Value = "https://www.betexplorer.com/football/costa-rica/primera-division/san-carlos-zeledon/t0FgCMWF/"
